{"js": "// Split two \"wall of text\" paragraphs (Programa / Bibliografia) into\n// lines separated by <w:br/> while keeping every segment inside a\n// single run, mirroring the target OOXML diff.\n\nfunction xmlEscape(s) {\n  return s\n    .replace(/&/g, \"&amp;\")\n    .replace(/</g, \"&lt;\")\n    .replace(/>/g, \"&gt;\");\n}\n\n// Builds the OOXML for a <w:p><w:r>...</w:r></w:p> paragraph whose run\n// alternates <w:t> segments with <w:br/> (no trailing break after the\n// last segment).\nfunction buildParagraphOoxml(segments) {\n  const runInner = segments\n    .map((seg, i) => {\n      const preserve = seg !== seg.trim() || seg.length === 0;\n      const spaceAttr = preserve ? ' xml:space=\"preserve\"' : \"\";\n      const t = `<w:t${spaceAttr}>${xmlEscape(seg)}</w:t>`;\n      return i < segments.length - 1 ? t + \"<w:br/>\" : t;\n    })\n    .join(\"\");\n\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    \"<pkg:xmlData>\" +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body><w:p><w:r>\" +\n    runInner +\n    \"</w:r></w:p></w:body>\" +\n    \"</w:document>\" +\n    \"</pkg:xmlData>\" +\n    \"</pkg:part>\" +\n    \"</pkg:package>\"\n  );\n}\n\nconst programaSegments = [\n  \"1. Sistemas de Informa\u00e7\u00e3o\",\n  \"1.1. Sistemas de Processamento de Informa\u00e7\u00f5es;\",\n  \"1.2. Sistemas de Informa\u00e7\u00f5es Gerenciais;\",\n  \"1.3. Sistema de Apoio \u00e0 Decis\u00e3o;\",\n  \"1.4. Sistemas de Informa\u00e7\u00e3o no Com\u00e9rcio Eletr\u00f4nico;\",\n  \"1.5. Sistemas de Informa\u00e7\u00e3o em Cadeia de Suprimentos;\",\n  \"1.6. Sistemas inteligentes nos neg\u00f3cios;\",\n  \"1.7. Sistemas estrat\u00e9gicos. \",\n  \"2. Projeto de Sistemas de Informa\u00e7\u00e3o.\",\n  \"2.1. Especifica\u00e7\u00e3o das Sa\u00eddas;\",\n  \"2.2. Especifica\u00e7\u00e3o dos Arquivos;\",\n  \"2.3. Especifica\u00e7\u00e3o das Entradas;\",\n  \"2.4. Especifica\u00e7\u00e3o do Processamento.\",\n  \"3. Tecnologia de Informa\u00e7\u00e3o.\",\n  \"3.1. Evolu\u00e7\u00e3o da Computa\u00e7\u00e3o;\",\n  \"3.2. Recursos Computacionais.\",\n  \"4. Processo de Desenvolvimento de Sistemas de Informa\u00e7\u00e3o.\",\n  \"4.1. Defini\u00e7\u00e3o do Neg\u00f3cio;\",\n  \"4.2. Identifica\u00e7\u00e3o do Problema e/ou Oportunidades;\",\n  \"4.3. Sele\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n  \"4.4. Implementa\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;\",\n  \"4.5. Avalia\u00e7\u00e3o da Efic\u00e1cia do Sistema de Informa\u00e7\u00e3o;\",\n];\n\nconst bibliografiaSegments = [\n  \"HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004\",\n  \"LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gest\u00e3o integrada de processos e da tecnologia da informa\u00e7\u00e3o. S\u00e3o Paulo:Atlas, 2006.\",\n  \"LAURINDO, F.J.B. Tecnologia da Informa\u00e7\u00e3o: Efic\u00e1cia nas Organiza\u00e7\u00f5es. S\u00e3o Paulo, Editora Futura, 2002.\",\n  \"STAIR, R.M., Princ\u00edpios de Sistema de Informa\u00e7\u00e3o: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.\",\n  \"TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.\",\n  \"TURBAN, E., RAIANER JR, K., POTTER, R. E., Administra\u00e7\u00e3o de Tecnologia da Informa\u00e7\u00e3o: Teoria e Pr\u00e1tica\u201d, S\u00e3o Paulo, Editora Campus, 2003.\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two target paragraphs by their (unique) original text so\n// the script does not depend on a brittle, hardcoded paragraph index.\nlet programaParagraph = null;\nlet bibliografiaParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (programaParagraph === null && text.indexOf(\"1. Sistemas de Informa\u00e7\u00e3o1.1.\") === 0) {\n    programaParagraph = para;\n  } else if (bibliografiaParagraph === null && text.indexOf(\"HAL R. VARIAN,\") === 0) {\n    bibliografiaParagraph = para;\n  }\n}\n\nif (!programaParagraph) {\n  throw new Error(\"Could not locate the 'Programa' paragraph.\");\n}\nif (!bibliografiaParagraph) {\n  throw new Error(\"Could not locate the 'Bibliografia' paragraph.\");\n}\n\nprogramaParagraph.insertOoxml(buildParagraphOoxml(programaSegments), Word.InsertLocation.replace);\nbibliografiaParagraph.insertOoxml(buildParagraphOoxml(bibliografiaSegments), Word.InsertLocation.replace);\n\nawait context.sync();\n", "ps1": "# Split two \"wall of text\" paragraphs (Programa / Bibliografia) into\n# lines separated by a line break (<w:br/>) while keeping every segment\n# inside a single run, mirroring the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n$programaXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>1. Sistemas de Informa\u00e7\u00e3o</w:t><w:br/><w:t>1.1. Sistemas de Processamento de Informa\u00e7\u00f5es;</w:t><w:br/><w:t>1.2. Sistemas de Informa\u00e7\u00f5es Gerenciais;</w:t><w:br/><w:t>1.3. Sistema de Apoio \u00e0 Decis\u00e3o;</w:t><w:br/><w:t>1.4. Sistemas de Informa\u00e7\u00e3o no Com\u00e9rcio Eletr\u00f4nico;</w:t><w:br/><w:t>1.5. Sistemas de Informa\u00e7\u00e3o em Cadeia de Suprimentos;</w:t><w:br/><w:t>1.6. Sistemas inteligentes nos neg\u00f3cios;</w:t><w:br/><w:t xml:space=\"preserve\">1.7. Sistemas estrat\u00e9gicos. </w:t><w:br/><w:t>2. Projeto de Sistemas de Informa\u00e7\u00e3o.</w:t><w:br/><w:t>2.1. Especifica\u00e7\u00e3o das Sa\u00eddas;</w:t><w:br/><w:t>2.2. Especifica\u00e7\u00e3o dos Arquivos;</w:t><w:br/><w:t>2.3. Especifica\u00e7\u00e3o das Entradas;</w:t><w:br/><w:t>2.4. Especifica\u00e7\u00e3o do Processamento.</w:t><w:br/><w:t>3. Tecnologia de Informa\u00e7\u00e3o.</w:t><w:br/><w:t>3.1. Evolu\u00e7\u00e3o da Computa\u00e7\u00e3o;</w:t><w:br/><w:t>3.2. Recursos Computacionais.</w:t><w:br/><w:t>4. Processo de Desenvolvimento de Sistemas de Informa\u00e7\u00e3o.</w:t><w:br/><w:t>4.1. Defini\u00e7\u00e3o do Neg\u00f3cio;</w:t><w:br/><w:t>4.2. Identifica\u00e7\u00e3o do Problema e/ou Oportunidades;</w:t><w:br/><w:t>4.3. Sele\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;</w:t><w:br/><w:t>4.4. Implementa\u00e7\u00e3o do Sistema de Informa\u00e7\u00e3o;</w:t><w:br/><w:t>4.5. Avalia\u00e7\u00e3o da Efic\u00e1cia do Sistema de Informa\u00e7\u00e3o;</w:t></w:r></w:p>'\n\n$bibliografiaXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r><w:t>HAL R. VARIAN, H. R.; FARRELL, J., SHAPIRO, C. The economics of information technology: an introduction. Cambridge University Press, 2004</w:t><w:br/><w:t>LAURINDO, F. J. B.; ROTONDARO, R. G. orgs. Gest\u00e3o integrada de processos e da tecnologia da informa\u00e7\u00e3o. S\u00e3o Paulo:Atlas, 2006.</w:t><w:br/><w:t>LAURINDO, F.J.B. Tecnologia da Informa\u00e7\u00e3o: Efic\u00e1cia nas Organiza\u00e7\u00f5es. S\u00e3o Paulo, Editora Futura, 2002.</w:t><w:br/><w:t>STAIR, R.M., Princ\u00edpios de Sistema de Informa\u00e7\u00e3o: Uma Abordagem Gerencial, Rio de Janeiro, LTC, 1998.</w:t><w:br/><w:t>TURBAN, E. et al. Information Technology for Management: Transforming Organizations in the Digital Economy. 7th edition, Wiley, 2009.</w:t><w:br/><w:t>TURBAN, E., RAIANER JR, K., POTTER, R. E., Administra\u00e7\u00e3o de Tecnologia da Informa\u00e7\u00e3o: Teoria e Pr\u00e1tica\u201d, S\u00e3o Paulo, Editora Campus, 2003.</w:t></w:r></w:p>'\n\n# Locate the two target paragraphs by their (unique) original text\n# prefix so the script does not depend on a brittle, hardcoded index.\n$programaIndex = 0\n$bibliografiaIndex = 0\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $t = $d.Paragraphs.Item($i).Range.Text\n  if ($programaIndex -eq 0 -and $t.StartsWith(\"1. Sistemas de Informa\u00e7\u00e3o1.1.\")) {\n    $programaIndex = $i\n  } elseif ($bibliografiaIndex -eq 0 -and $t.StartsWith(\"HAL R. VARIAN,\")) {\n    $bibliografiaIndex = $i\n  }\n}\n\nif ($programaIndex -eq 0) {\n  throw \"Could not locate the 'Programa' paragraph.\"\n}\nif ($bibliografiaIndex -eq 0) {\n  throw \"Could not locate the 'Bibliografia' paragraph.\"\n}\n\n[void]$d.Paragraphs.Item($programaIndex).Range.InsertXML($programaXml)\n[void]$d.Paragraphs.Item($bibliografiaIndex).Range.InsertXML($bibliografiaXml)\n"}
